$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row below the header row, pushing the existing metadata rows
# down by one. This makes room for a new row of "machine name" identifiers
# that relate each column's human label (row 1) to a short code used to
# build SKOS hierarchical relations between columns.
$ws.Rows("2:2").Insert()

# The row that used to be a stray, mostly-empty row 5 (only H5 populated
# with "mapping-ano.xlsx") shifted down to row 6 when we inserted above;
# remove it since its real content now lives in the (shifted) row 5.
$ws.Rows("6:6").Delete()

# Populate the newly inserted row 2 with the machine-readable column names.
$ws.Range("A2").Value = "subseccion-codigo"
$ws.Range("B2").Value = "comarca-nombre"
$ws.Range("C2").Value = "vab"
$ws.Range("D2").Value = "sector-vab-descripcion"
$ws.Range("E2").Value = "subseccion-descripcion"
$ws.Range("F2").Value = "comarca-codigo"
$ws.Range("G2").Value = "sector-vab-codigo"
$ws.Range("H2").Value = "ano"
